$d = $word.ActiveDocument

# Locate the bold "Java" run at the end of the "Basic knowledge of ..." line
# (the only occurrence of the whole word "Java" in the document).
$rng = $d.Content
$found = $rng.Find.Execute("Java", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'Java' in the document"
}

# Collapse to the end of the found text so we can append after it.
$rng.Collapse(0)

# ","  (plain, not bold)
$rng.InsertAfter(",")
$rng.Collapse(0)

# " "  (plain, not bold)
$rng.InsertAfter(" ")
$rng.Collapse(0)

# "HTML" (bold) - set Bold right after the insert, while $rng still spans
# the freshly-inserted text (setting Bold on a collapsed/empty range is
# unsafe and can affect unrelated text).
$rng.InsertAfter("HTML")
$rng.Bold = 1
$rng.Collapse(0)

# ","  (plain, not bold)
$rng.InsertAfter(",")
$rng.Collapse(0)

# " "  (plain, not bold)
$rng.InsertAfter(" ")
$rng.Collapse(0)

# "CSS" (bold)
$rng.InsertAfter("CSS")
$rng.Bold = 1
$rng.Collapse(0)
